$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ25005660",
    "summ25104368",
    "summ25209694",
    "summ25319269",
    "summ25451353",
    "summ25583380",
    "summ25728628",
    "summ25873621",
    "summ26011790",
    "summ26146306",
    "summ26289626",
    "summ26434507",
    "summ26574957",
    "summ26713434",
    "summ26848945",
    "summ26986645",
    "summ27129160",
    "summ27271769",
    "summ27414214",
    "summ27589144",
    "summ27757712",
    "summ27918430",
    "summ28048934",
    "summ28193003",
    "summ28326981",
    "summ28462716",
    "summ28603050",
    "summ28748332",
    "summ28893603",
    "summ29033112",
    "summ29177704",
    "summ29317448",
    "summ29490824",
    "summ29673365",
    "summ29829818",
    "summ30006746",
    "summ30158710",
    "summ30303695",
    "summ30446208",
    "summ30601683",
    "summ30742333",
    "summ30888016",
    "summ31026063",
    "summ31162764",
    "summ31315249",
    "summ31462607",
    "summ31594835",
    "summ31737994",
    "summ31875081",
    "summ32018646"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
